# Updated cryptos list (coin prices / 1h volume deltas refreshed by the
# scraping Action). Most cells are plain text replacements; a handful of
# "Price" column values are numeric-looking strings ("1.00", "605.76", ...)
# that Excel's COM layer would otherwise silently coerce into real numbers
# (dropping trailing zeros / adding float noise) when assigned through
# .Value. For those we briefly force a text NumberFormat while writing the
# value, then restore the cell's default ("Normal") style so no visible
# formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.064.16"
$ws.Range("E2").Value = "  +1.03%  "

$ws.Range("D3").Value = "2.747.57"
$ws.Range("E3").Value = "  +3.60%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.98%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  +0.58%  "

$ws.Range("D9").Value = "2.748.01"
$ws.Range("E9").Value = "  +3.60%  "

$ws.Range("E10").Value = "  -1.73%  "

# Rows 11/12 swap ranking order: Toncoin now above Cardano.
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.37"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.01%  "

$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.365"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.75%  "

$ws.Range("E13").Value = "  -0.38%  "

$ws.Range("E14").Value = "  +2.95%  "

$ws.Range("D15").Value = "3.248.98"
$ws.Range("E15").Value = "  +3.65%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000190"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.06%  "

$ws.Range("D17").Value = "69.009.01"
$ws.Range("E17").Value = "  +1.08%  "

$ws.Range("D18").Value = "2.746.83"
$ws.Range("E18").Value = "  +3.43%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.03%  "

# Rows 20/21 swap ranking order: BitcoinCash now above Uniswap.
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "369.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.35%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.35%  "

$ws.Range("E22").Value = "  +3.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.11%  "

$ws.Range("E24").Value = "  +3.43%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.26%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("E27").Value = "  +2.70%  "

$ws.Range("E29").Value = "  +1.61%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "602.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.72%  "

$ws.Range("E31").Value = "  -2.62%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.92%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.88%  "

$ws.Range("E34").Value = "  +6.07%  "

$ws.Range("E35").Value = "  +3.33%  "

$ws.Range("E36").Value = "  +4.47%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "20.21"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "163.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.43%  "

$ws.Range("E40").Value = "  +3.69%  "

$ws.Range("E41").Value = "  +2.49%  "

$ws.Range("E42").Value = "  +2.99%  "

$ws.Range("E43").Value = "  +3.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.28%  "

$ws.Range("E45").Value = "  -4.97%  "

$ws.Range("E46").Value = "  +0.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "159.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.73%  "

$ws.Range("E48").Value = "  +5.31%  "

$ws.Range("E49").Value = "  +6.83%  "

$ws.Range("E50").Value = "  +7.79%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.47%  "
